$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Columns whose values actually differ between row 20 and row 21.
# (Other columns, e.g. D, I, T, U, V, W, Y, AA, AD, AE, AG, AT, AY, are
# identical between the two rows and do not need to be touched.)
$cols = @("A","B","E","F","G","H","M","P","Q","R","S","Z","AB","AC","AW","AX")

foreach ($col in $cols) {
    $c20 = $ws.Range($col + "20")
    $c21 = $ws.Range($col + "21")

    $v20 = $c20.Value2
    $v21 = $c21.Value2

    $c20.Value2 = $v21
    $c21.Value2 = $v20
}
